$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: value edits on rows 2-25 (no row shifting needed for these) ---
$ws.Range("C2").Value = 14.9
$ws.Range("F3").ClearContents()
$ws.Range("F4").Value = 17.97
$ws.Range("F5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("C12").Value = 12.5
$ws.Range("C14").ClearContents()
$ws.Range("C20").Value = 12.5
$ws.Range("C21").Value = 12.7
$ws.Range("C22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("F23").Value = 16.48

# --- Step 2: remove the "RM 232" row (row 26) and the "SC 92" row (which
# becomes row 27 once row 26 is removed) entirely, shifting later rows up ---
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- Step 3: value edits on the now-shifted rows 27-33 ---
$ws.Range("F27").ClearContents()
$ws.Range("F29").Value = 18.06
$ws.Range("B30").Value = -19.7
$ws.Range("C31").Value = 15.3
$ws.Range("B32").ClearContents()
$ws.Range("C33").Value = 10.4
